# Edit script: applies the changes described by the diff
# 1. Rename B1 headers on "Weekly Quantity" and "Monthly Trend" sheets
# 2. Add a new "PO Forecast" worksheet with forecast data

$wb = $excel.ActiveWorkbook

# --- Step 1: update header text on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Step 2: add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$headerVals = New-Object 'object[,]' 1,4
$headerVals[0,0] = "ds"
$headerVals[0,1] = "PO_Forecast"
$headerVals[0,2] = "yhat_lower"
$headerVals[0,3] = "yhat_upper"
$wsForecast.Range("A1:D1").Value = $headerVals

$data = New-Object 'object[,]' 31,4
$data[0,0] = 45410.99999999999
$data[0,1] = 15
$data[0,2] = 1.44711643166707
$data[0,3] = 28.25071143200483
$data[1,0] = 45417.99999999999
$data[1,1] = 15
$data[1,2] = 2.161985867679756
$data[1,3] = 27.47628161171109
$data[2,0] = 45424.99999999999
$data[2,1] = 15
$data[2,2] = 0.9795329760564324
$data[2,3] = 27.25142136285471
$data[3,0] = 45431.99999999999
$data[3,1] = 15
$data[3,2] = 1.619333466015128
$data[3,3] = 28.25832050999749
$data[4,0] = 45438.99999999999
$data[4,1] = 15
$data[4,2] = 1.630903898728106
$data[4,3] = 28.17124480221236
$data[5,0] = 45445.99999999999
$data[5,1] = 15
$data[5,2] = 1.711852971308123
$data[5,3] = 27.28884273761127
$data[6,0] = 45452.99999999999
$data[6,1] = 14
$data[6,2] = 1.041049515054502
$data[6,3] = 28.03413768939986
$data[7,0] = 45459.99999999999
$data[7,1] = 14
$data[7,2] = 1.19126312647168
$data[7,3] = 27.07143172697685
$data[8,0] = 45466.99999999999
$data[8,1] = 14
$data[8,2] = 0.3904633536369073
$data[8,3] = 26.95124570993538
$data[9,0] = 45480.99999999999
$data[9,1] = 14
$data[9,2] = -0.5000196530637044
$data[9,3] = 27.12004911808427
$data[10,0] = 45487.99999999999
$data[10,1] = 14
$data[10,2] = 1.281908590795265
$data[10,3] = 26.18802679706323
$data[11,0] = 45522.99999999999
$data[11,1] = 13
$data[11,2] = 0.3094734141033493
$data[11,3] = 25.07938236604395
$data[12,0] = 45529.99999999999
$data[12,1] = 13
$data[12,2] = -0.2679262945930946
$data[12,3] = 26.30299706157736
$data[13,0] = 45536.99999999999
$data[13,1] = 13
$data[13,2] = 0.4777557647614131
$data[13,3] = 25.27301331674511
$data[14,0] = 45550.99999999999
$data[14,1] = 12
$data[14,2] = 0.05262064230914584
$data[14,3] = 24.71230521850417
$data[15,0] = 45557.99999999999
$data[15,1] = 12
$data[15,2] = -0.672162357668087
$data[15,3] = 25.17202318438308
$data[16,0] = 45564.99999999999
$data[16,1] = 12
$data[16,2] = -1.056131629905958
$data[16,3] = 25.69292224063442
$data[17,0] = 45571.99999999999
$data[17,1] = 12
$data[17,2] = -0.3312858482844827
$data[17,3] = 25.24719458235435
$data[18,0] = 45578.99999999999
$data[18,1] = 12
$data[18,2] = -1.611085889392447
$data[18,3] = 24.2782643002923
$data[19,0] = 45599.99999999999
$data[19,1] = 11
$data[19,2] = -0.9823454525308103
$data[19,3] = 23.95093471424883
$data[20,0] = 45613.99999999999
$data[20,1] = 11
$data[20,2] = -2.798168283215304
$data[20,3] = 23.98239070122738
$data[21,0] = 45620.99999999999
$data[21,1] = 11
$data[21,2] = -1.615905096427415
$data[21,3] = 23.60116703428174
$data[22,0] = 45627.99999999999
$data[22,1] = 11
$data[22,2] = -1.955682878193036
$data[22,3] = 23.83827608680078
$data[23,0] = 45634.99999999999
$data[23,1] = 11
$data[23,2] = -2.789381178041031
$data[23,3] = 22.66584761826051
$data[24,0] = 45641.99999999999
$data[24,1] = 10
$data[24,2] = -2.518413885324669
$data[24,3] = 22.47297183385867
$data[25,0] = 45648.99999999999
$data[25,1] = 10
$data[25,2] = -3.289439638468262
$data[25,3] = 23.48552350283035
$data[26,0] = 45655.99999999999
$data[26,1] = 10
$data[26,2] = -2.899399363522755
$data[26,3] = 22.56611059041239
$data[27,0] = 45662.99999999999
$data[27,1] = 10
$data[27,2] = -2.508598386166035
$data[27,3] = 23.16917894300598
$data[28,0] = 45669.99999999999
$data[28,1] = 10
$data[28,2] = -3.114923774791548
$data[28,3] = 21.91103771148572
$data[29,0] = 45676.99999999999
$data[29,1] = 10
$data[29,2] = -2.6290878803665
$data[29,3] = 22.92023966146159
$data[30,0] = 45683.99999999999
$data[30,1] = 9
$data[30,2] = -3.526660922023199
$data[30,3] = 21.60537192875155

$wsForecast.Range("A2:D32").Value = $data

# --- Step 3: copy cell formatting from the "Weekly Quantity" sheet so the
#     new sheet matches the house style (bold header w/ border, date format) ---
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A32").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

Write-Output "PO Forecast sheet created and headers updated."
